$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("ND-SECOND-YEAR-SECOND-SEMESTER")
$ws4.Range("A8").Value = "EED216"

$ws2 = $wb.Worksheets.Item("ND-FIRST-YEAR-SECOND-SEMESTER")
$ws2.Range("A9").Value = "EED126"

$ws4.Select()
$ws4.Range("A8").Select()

$ws2.Select()
$ws2.Range("E8").Select()

$wb.Save()
